$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (reflects new "through" date)
$ws.Name = "Through 2021-09-18"

# Update the label for September row
$ws.Range("A10").Value = "September (through 09-18)"

# Update September row (row 10) values
$ws.Range("B10").Value = 20
$ws.Range("C10").Value = 33
$ws.Range("D10").Value = 41
$ws.Range("E10").Value = 34
$ws.Range("F10").Value = 41
$ws.Range("G10").Value = 66
$ws.Range("H10").Value = 95

# Update Total row (row 11) values
$ws.Range("B11").Value = 214
$ws.Range("C11").Value = 414
$ws.Range("D11").Value = 592
$ws.Range("E11").Value = 524
$ws.Range("F11").Value = 390
$ws.Range("G11").Value = 850
$ws.Range("H11").Value = 1165
